# Generate Report for Handoff
#
# A new handoff report was generated, which refreshes the "Latest Handoff
# Datetime" for the file that is currently "Ready for handoff" and whose
# handoff package was just (re)generated:
#   48430eff-6746-4dfc-b2d7-cbb467fa8e4c.md (row 5 on both the zh-cn and
#   de-de localization-status sheets).
#
# zh-cn : E5  2016-03-20 22:35:18  ->  2016-03-20 22:35:37
# de-de : E5  2016-03-20 22:35:21  ->  2016-03-20 22:35:40

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E5").Value = "2016-03-20 22:35:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E5").Value = "2016-03-20 22:35:40"
